$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values: the message recipient email(s) changed
$ws.Range("B1").Value = "capin19887@kernuo.com"
$ws.Range("B2").Value = "komajox883@fna6.com"

# Move the selection to B2, reflecting the single-user send target
$ws.Range("B2").Select()
